# Update "想去人数" (F) and "最低票价" (G) figures across the four sheets
# of the 广州-漫展信息 workbook to match the newly scraped data snapshot.

$wb = $excel.ActiveWorkbook

# ---- Sheet "展览" (exhibitions) ----
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1872
$ws1.Range("G2").Value = 58
$ws1.Range("F3").Value = 1508
$ws1.Range("F4").Value = 875
$ws1.Range("F5").Value = 761
$ws1.Range("F6").Value = 13293
$ws1.Range("F7").Value = 13169
$ws1.Range("F8").Value = 1014
$ws1.Range("F9").Value = 772
$ws1.Range("F13").Value = 667
$ws1.Range("F14").Value = 2089
$ws1.Range("F17").Value = 68
$ws1.Range("F18").Value = 19
$ws1.Range("F19").Value = 392
$ws1.Range("F20").Value = 240
$ws1.Range("F21").Value = 284
$ws1.Range("F22").Value = 414
$ws1.Range("F23").Value = 752
$ws1.Range("F24").Value = 13

# ---- Sheet "演出" (performances) ----
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 18
$ws2.Range("F5").Value = 130
$ws2.Range("F6").Value = 60
$ws2.Range("F7").Value = 117
$ws2.Range("F9").Value = 25

# ---- Sheet "本地生活" (local life) ----
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F3").Value = 34

# ---- Sheet "全部类型" (all types, aggregated view) ----
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 1872
$ws4.Range("G3").Value = 58
$ws4.Range("F4").Value = 1508
$ws4.Range("F5").Value = 875
$ws4.Range("F7").Value = 761
$ws4.Range("F8").Value = 13293
$ws4.Range("F9").Value = 13169
$ws4.Range("F10").Value = 1014
$ws4.Range("F11").Value = 772
$ws4.Range("F15").Value = 667
$ws4.Range("F17").Value = 18
$ws4.Range("F18").Value = 2089
$ws4.Range("F21").Value = 68
$ws4.Range("F22").Value = 130
$ws4.Range("F23").Value = 19
$ws4.Range("F24").Value = 60
$ws4.Range("F25").Value = 34
$ws4.Range("F26").Value = 392
$ws4.Range("F27").Value = 240
$ws4.Range("F28").Value = 284
$ws4.Range("F29").Value = 414
$ws4.Range("F30").Value = 752
$ws4.Range("F31").Value = 117
$ws4.Range("F33").Value = 13
$ws4.Range("F34").Value = 25

$wb.Save()
